$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meta")

$ws.Range("C1").Value = "RootListPath iati-organisation"

$ws.Range("C2").Select()
